$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 13:46"

# Refresh Canary Islands (and a few other provinces) case counts, which changes the
# descending sort order by "Casos totales" (column B) for the affected rows.
# Row 25
$ws.Range("A25").Value = "Tenerife"
$ws.Range("B25").Value = 293
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 284
$ws.Range("E25").Value = 11
# Row 26
$ws.Range("A26").Value = "Burgos"
$ws.Range("B26").Value = 269
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 175
$ws.Range("E26").Value = 16
# Row 27
$ws.Range("A27").Value = "Salamanca"
$ws.Range("B27").Value = 265
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 180
$ws.Range("E27").Value = 21
# Row 28
$ws.Range("A28").Value = "Guadalajara"
$ws.Range("B28").Value = 263
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 257
$ws.Range("E28").Value = 4
# Row 39
$ws.Range("A39").Value = "Gran Canaria"
$ws.Range("B39").Value = 135
$ws.Range("C39").Value = 7
$ws.Range("D39").Value = 133
$ws.Range("E39").Value = 11
# Row 40
$ws.Range("A40").Value = "Cadiz"
$ws.Range("B40").Value = 134
$ws.Range("C40").Value = 72
$ws.Range("D40").Value = 126
$ws.Range("E40").Value = 2
# Row 41
$ws.Range("A41").Value = "Cuenca"
$ws.Range("B41").Value = 120
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = 104
$ws.Range("E41").Value = 8
# Row 54
$ws.Range("A54").Value = "La Palma"
$ws.Range("C54").Value = 7
$ws.Range("D54").Value = 21
$ws.Range("E54").Value = 11
# Row 55
$ws.Range("A55").Value = "Ibiza"
$ws.Range("B55").Value = 21
$ws.Range("D55").Value = 20
$ws.Range("E55").Value = 1
# Row 56
$ws.Range("B56").Value = 18
$ws.Range("C56").Value = 7
$ws.Range("D56").Value = 18
# Row 57
$ws.Range("A57").Value = "Menorca"
$ws.Range("B57").Value = 15
$ws.Range("C57").Value = 10
$ws.Range("D57").Value = 13
$ws.Range("E57").Value = 0
# Row 58
$ws.Range("C58").Value = 7
$ws.Range("D58").Value = 9
# Row 61
$ws.Range("C61").Value = 7
# Row 62
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 2
